# Add team record (Wins/Losses/Ties) columns AD, AE, AF to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from an existing header cell (AC1) onto the
# new header cells so they match the look of the rest of the header row
# (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record is the same for every player row (2-50): 82 wins, 79 losses, 0 ties.
$lastRow = 50
$winsRange = $ws.Range("AD2:AD" + $lastRow)
$lossesRange = $ws.Range("AE2:AE" + $lastRow)
$tiesRange = $ws.Range("AF2:AF" + $lastRow)

$winsRange.Value = 82
$lossesRange.Value = 79
$tiesRange.Value = 0
